# Manual testing of Experiment-Methodology-04 and aggregation of the Consolidation Theory
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Sheet1: add the four new "Experiment" data rows (25-28) plus the totals
#    row (29) below the existing data.
# ---------------------------------------------------------------------------
$ws1.Range("I25").Value = 962
$ws1.Range("J25").Value = 19232
$ws1.Range("K25").Formula = "=J25/200000"

$ws1.Range("I26").Value = 2497
$ws1.Range("J26").Value = 21988
$ws1.Range("K26").Formula = "=J26/200000"

$ws1.Range("I27").Value = 1952
$ws1.Range("J27").Value = 23941
$ws1.Range("K27").Formula = "=J27/200000"

$ws1.Range("I28").Value = 2005
$ws1.Range("J28").Value = 18088
$ws1.Range("K28").Formula = "=J28/200000"

$ws1.Range("I29").Formula = "=SUM(I25:I28)"
$ws1.Range("J29").Formula = "=SUM(J25:J28)"
$ws1.Range("K29").Formula = "=SUM(K25:K28)"

# Column width changes on Sheet1 (col B widened, new cols F/G sized, cols C:E
# keep their old width).
$ws1.Columns.Item(2).ColumnWidth = 40.5
$ws1.Columns.Item(6).ColumnWidth = 9.8
$ws1.Columns.Item(7).ColumnWidth = 12.5

# ---------------------------------------------------------------------------
# 2. Add the new worksheet "Experiment-Methodology-04" right after Sheet1.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Experiment-Methodology-04"

# NOTE: the string literals below are written in the precise order needed to
# reproduce the author's original shared-string table ordering (new strings
# are appended to xl/sharedStrings.xml in first-use order). Numeric/formula
# cells can be written in any order since they don't touch that table.

# -- headers (first-use order establishes shared-string indices 22-24) ------
$ws2.Range("B2").Value = "Pre-load Files"
$ws2.Range("B12").Value = "Dynamic Files"
$ws2.Range("C2").Value = "Token Count"
$ws2.Range("D2").Value = "Lines"
$ws2.Range("E2").Value = "Total Tokens"
$ws2.Range("B6").Value = "Pre-load Files"
$ws2.Range("C6").Value = "Token Count"
$ws2.Range("D6").Value = "Lines"
$ws2.Range("E6").Value = "Total Tokens"
$ws2.Range("C12").Value = "Token Count"
$ws2.Range("D12").Value = "Lines"
$ws2.Range("E12").Value = "Total Tokens"

# -- Table 2 (specs) file names: indices 25-28 --------------------------------
$ws2.Range("B7").Value = "docs/specs/operations-manual-standard.md"
$ws2.Range("B8").Value = "docs/specs/operations-manual-exceptions.md"
$ws2.Range("B9").Value = "docs/specs/architecture-deep-dive.md"
$ws2.Range("B10").Value = "docs/specs/troubleshooting-compendium.md"

# -- Table 3 (dynamic files) file names: indices 29-35 -------------------------
$ws2.Range("B13").Value = "Read(docs/wpds/pipeline-refactor.md)"
$ws2.Range("B14").Value = "Read(docs/specs/data-pipeline-overview.md)"
$ws2.Range("B15").Value = "Read(docs/specs/module-alpha.md)"
$ws2.Range("B16").Value = "Read(docs/specs/module-beta.md)"
$ws2.Range("B17").Value = "Read(docs/specs/module-gamma.md)"
$ws2.Range("B18").Value = "Read(docs/specs/integration-layer.md)"
$ws2.Range("B19").Value = "Read(docs/specs/compliance-requirements.md)"

# -- remaining headers: indices 36-37 ------------------------------------------
$ws2.Range("F2").Value = "size (bytes)"
$ws2.Range("G2").Value = "byte-to-tokens"
$ws2.Range("F6").Value = "size (bytes)"
$ws2.Range("G6").Value = "byte-to-tokens"
$ws2.Range("F12").Value = "size (bytes)"
$ws2.Range("G12").Value = "byte-to-tokens"

$ws2.Range("B2:G2").Font.Bold = $true
$ws2.Range("B6:G6").Font.Bold = $true
$ws2.Range("B12:G12").Font.Bold = $true

# -- token/byte helper labels: indices 38-40 -----------------------------------
$ws2.Range("G23").Value = "Tokens -> "
$ws2.Range("H23").Value = "Bytes"
$ws2.Range("G25").Value = "Bytes -> "
$ws2.Range("H25").Value = "Tokens"
$ws2.Range("G23:H23").Font.Bold = $true
$ws2.Range("G25:H25").Font.Bold = $true

# -- Totals label: index 41 -----------------------------------------------------
$ws2.Range("B21").Value = "Totals"

# -- Table 1 file names (entered last): indices 42-43 ---------------------------
$ws2.Range("B3").Value = ".claude/commands/setup-hard.md"
$ws2.Range("B4").Value = ".claude/commands/analyze-wpd.md"

# ---------------------------------------------------------------------------
# Numeric values and formulas (order irrelevant to shared-string table).
# ---------------------------------------------------------------------------
# --- Table 1: Pre-load Files (local commands) ---------------------------------
$ws2.Range("C3").Value = 402
$ws2.Range("D3").Value = 52
$ws2.Range("F3").Value = 1490
$ws2.Range("G3").Formula = "=F3/C3"

$ws2.Range("C4").Value = 618
$ws2.Range("D4").Value = 63
$ws2.Range("E4").Formula = "=SUM(C3:C4)"
$ws2.Range("F4").Value = 2559
$ws2.Range("G4").Formula = "=F4/C4"

# --- Table 2: Pre-load Files (specs) ------------------------------------------
$ws2.Range("C7").Value = 19323
$ws2.Range("D7").Value = 963
$ws2.Range("F7").Value = 108497
$ws2.Range("G7").Formula = "=F7/C7"

$ws2.Range("C8").Value = 15636
$ws2.Range("D8").Value = 1593
$ws2.Range("E8").Formula = "=SUM(C7:C8)"
$ws2.Range("F8").Value = 66444
$ws2.Range("G8").Formula = "=F8/C8"

$ws2.Range("C9").Value = 14676
$ws2.Range("D9").Value = 1071
$ws2.Range("E9").Formula = "=SUM(C7:C9)"
$ws2.Range("F9").Value = 85873
$ws2.Range("G9").Formula = "=F9/C9"

$ws2.Range("C10").Value = 18477
$ws2.Range("D10").Value = 2006
$ws2.Range("E10").Formula = "=SUM(C7:C10)"
$ws2.Range("F10").Value = 83359
$ws2.Range("G10").Formula = "=F10/C10"

# --- Table 3: Dynamic Files ----------------------------------------------------
$ws2.Range("C13").Value = 5034
$ws2.Range("D13").Value = 393
$ws2.Range("F13").Value = 21978
$ws2.Range("G13").Formula = "=F13/C13"

$ws2.Range("C14").Value = 6041
$ws2.Range("D14").Value = 426
$ws2.Range("F14").Value = 32350
$ws2.Range("G14").Formula = "=F14/C14"

$ws2.Range("C15").Value = 6204
$ws2.Range("D15").Value = 743
$ws2.Range("F15").Value = 24863
$ws2.Range("G15").Formula = "=F15/C15"

$ws2.Range("C16").Value = 6198
$ws2.Range("D16").Value = 742
$ws2.Range("F16").Value = 26029
$ws2.Range("G16").Formula = "=F16/C16"

$ws2.Range("C17").Value = 7658
$ws2.Range("D17").Value = 772
$ws2.Range("F17").Value = 33133
$ws2.Range("G17").Formula = "=F17/C17"

$ws2.Range("C18").Value = 4886
$ws2.Range("D18").Value = 531
$ws2.Range("F18").Value = 20349
$ws2.Range("G18").Formula = "=F18/C18"

$ws2.Range("C19").Value = 3939
$ws2.Range("D19").Value = 393
$ws2.Range("E19").Formula = "=SUM(C13:C19)"
$ws2.Range("F19").Value = 18690
$ws2.Range("G19").Formula = "=F19/C19"

# --- Totals row -----------------------------------------------------------------
$ws2.Range("C21").Formula = "=SUM(C7:C10,C13:C19)"
$ws2.Range("D21").Formula = "=SUM(D7:D10,D13:D19)"
$ws2.Range("F21").Formula = "=SUM(F7:F10,F13:F19)"
$ws2.Range("G21").Formula = "=F21/C21"

# --- Token / byte conversion helper table ---------------------------------------
$ws2.Range("G24").Value = 25000
$ws2.Range("G24").Interior.Color = 65535
$ws2.Range("H24").Formula = "=G24*G21"
$ws2.Range("H24").NumberFormat = "0"

$ws2.Range("G26").Value = 66000
$ws2.Range("G26").Interior.Color = 65535
$ws2.Range("H26").Formula = "=G26/G21"
$ws2.Range("H26").NumberFormat = "0"

# Column widths for the new sheet.
$ws2.Columns.Item(2).ColumnWidth = 40.5
$ws2.Range("C2:G19").ColumnWidth = 13
$ws2.Columns.Item(8).ColumnWidth = 12.17

# Selection / active cell state to match the final view.
$ws1.Range("J33").Select()
$ws2.Range("B25").Select()

# Window geometry (best effort; matches the author's recorded window state).
$excel.ActiveWindow.Left = 79780
$excel.ActiveWindow.Top = 6260
$excel.ActiveWindow.Width = 34380
$excel.ActiveWindow.Height = 21800
